$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 'section', 'Preparation and Environment '),
    @(3, 'protonation method', 'PROPKA'),
    @(4, 'pH', '7.4'),
    @(5, 'step type', 'conditional'),
    @(6, 'flow type', 'if'),
    @(7, 'flow parameter', 'membrane simulation'),
    @(8, 'flow logical parameter', 'e'),
    @(9, 'flow compared value', 'true'),
    @(10, 'Lipid type', 'POPC'),
    @(11, 'step type', 'conditional'),
    @(12, 'flow type', 'if'),
    @(13, 'flow parameter', 'membrane simulation'),
    @(14, 'flow logical parameter', 'e'),
    @(15, 'flow compared value', 'true'),
    @(16, 'box type', 'rectangular'),
    @(17, 'step type', 'conditional'),
    @(18, 'flow type', 'if'),
    @(19, 'flow parameter', 'membrane simulation'),
    @(20, 'flow logical parameter', 'e'),
    @(21, 'flow compared value', 'false'),
    @(22, 'box type', 'octahedral'),
    @(23, 'water type', 'TIP3P'),
    @(24, 'shell radius', '12 Å'),
    @(25, 'simulation', 'molecular dynamics'),
    @(26, 'suite', 'AMBER14'),
    @(27, 'step type', 'conditional'),
    @(28, 'flow type', 'if'),
    @(29, 'flow parameter', 'water type'),
    @(30, 'flow logical parameter', 'e'),
    @(31, 'flow compared value', 'TIP3P'),
    @(32, 'force field', 'ff14SB'),
    @(33, 'step type', 'conditional'),
    @(34, 'flow type', 'ELIF'),
    @(35, 'flow parameter', 'water type'),
    @(36, 'flow logical parameter', 'e'),
    @(37, 'flow compared value', ' OPC'),
    @(38, 'force field', 'ff19SB'),
    @(39, 'step type', 'conditional'),
    @(40, 'flow type', 'if'),
    @(41, 'flow parameter', 'membrane simulation'),
    @(42, 'flow logical parameter', 'e'),
    @(43, 'flow compared value', 'true'),
    @(44, 'lipid force field', 'LIPID14'),
    @(45, 'step type', 'conditional'),
    @(46, 'flow type', 'if'),
    @(47, 'flow parameter', 'membrane simulation'),
    @(48, 'flow logical parameter', 'e'),
    @(49, 'flow compared value', 'true'),
    @(50, 'dt', '2 fs'),
    @(51, 'cut', '9 Å'),
    @(52, 'step type', 'conditional'),
    @(53, 'flow type', 'if'),
    @(54, 'flow parameter', 'membrane simulation'),
    @(55, 'flow logical parameter', 'e'),
    @(56, 'flow compared value', 'false'),
    @(57, 'dt', '4 fs'),
    @(58, 'cut', '8 Å'),
    @(59, 'section', 'Minimization'),
    @(60, 'maxcyc', '17,500'),
    @(61, 'step type', 'iteration'),
    @(62, 'flow type', 'for each'),
    @(63, 'flow parameter', 'cycles of minimization print'),
    @(64, 'maxcyc', '2500'),
    @(65, 'restraint_wt', '25 kcal mol-1 Å-2'),
    @(66, 'restrainmask', 'MANUAL_INPUT'),
    @(67, 'section', 'Thermalization'),
    @(68, 'simulation time', '50'),
    @(69, 'nstlim', '12,500'),
    @(70, 'dt', '4 fs'),
    @(71, 'step type', 'conditional'),
    @(72, 'flow type', 'if'),
    @(73, 'flow parameter', 'ntp'),
    @(74, 'flow logical parameter', 'e'),
    @(75, 'flow compared value', '0'),
    @(76, 'MD', 'NVT'),
    @(77, 'ntp', '0'),
    @(78, 'step type', 'conditional'),
    @(79, 'flow type', 'elif'),
    @(80, 'flow parameter', 'ntp'),
    @(81, 'flow logical parameter', 'gt'),
    @(82, 'flow compared value', '0'),
    @(83, 'MD', 'NPT'),
    @(84, 'temp0', '100 K'),
    @(85, 'simulation time', '300'),
    @(86, 'nstlim', '75,000'),
    @(87, 'dt', '4 fs'),
    @(88, 'pres0', '1 atm'),
    @(89, 'temp0', '300 K'),
    @(90, 'restraint_wt', '10 kcal mol-1 Å-2'),
    @(91, 'restrainmask', 'MANUAL_INPUT'),
    @(92, 'simulation time', '300'),
    @(93, 'nstlim', '75,000'),
    @(94, 'dt', '4 fs'),
    @(95, 'step type', 'conditional'),
    @(96, 'flow type', 'if'),
    @(97, 'flow parameter', 'ntp'),
    @(98, 'flow logical parameter', 'e'),
    @(99, 'flow compared value', '0'),
    @(100, 'MD', 'NVT'),
    @(101, 'ntp', '0'),
    @(102, 'step type', 'conditional'),
    @(103, 'flow type', 'elif'),
    @(104, 'flow parameter', 'ntp'),
    @(105, 'flow logical parameter', 'gt'),
    @(106, 'flow compared value', '0'),
    @(107, 'MD', 'NPT'),
    @(108, 'restrainmask', 'MANUAL_INPUT'),
    @(109, 'restraint_wt', '0 kcal mol-1 Å-2'),
    @(110, 'simulation time', '100'),
    @(111, 'nstlim', '25,000'),
    @(112, 'dt', '4 fs'),
    @(113, 'section', 'Production'),
    @(114, 'overall repetitions', '5'),
    @(115, 'step type', 'conditional'),
    @(116, 'flow type', 'if'),
    @(117, 'flow parameter', 'ntp'),
    @(118, 'flow logical parameter', 'e'),
    @(119, 'flow compared value', '0'),
    @(120, 'MD', 'NVT'),
    @(121, 'ntp', '0'),
    @(122, 'step type', 'conditional'),
    @(123, 'flow type', 'elif'),
    @(124, 'flow parameter', 'ntp'),
    @(125, 'flow logical parameter', 'gt'),
    @(126, 'flow compared value', '0'),
    @(127, 'MD', 'NPT'),
    @(128, 'simulation time', '2'),
    @(129, 'nstlim', '500,000'),
    @(130, 'dt', '4 fs')
)

foreach ($row in $data) {
    $lineRow = $row[0]
    $ws.Cells.Item($lineRow, 2).Value = $row[1]
    $ws.Cells.Item($lineRow, 3).Value = $row[2]
}
